$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting (styles only) of the last existing week-block
#     (rows 48:55 = blank separator + 7-row block) down onto the new
#     block (rows 56:63) so fonts/borders/fills match the sheet's pattern.
$ws.Range("A48:D55").Copy()
$ws.Range("A56:D63").PasteSpecial(-4122)

# --- Blank separator row 56 keeps its style only (no values, like rows
#     8/16/24/32/40/48).

# --- New week header (row 57): date + column headers
$ws.Range("A57").Value = "04.11.17"
$ws.Range("B57").Value = "Calvin"
$ws.Range("C57").Value = "Finn"
$ws.Range("D57").Value = "Madi"

# --- Category rows 58-63
$ws.Range("A58").Value = "Konzept"
$ws.Range("A59").Value = "Meetings"
$ws.Range("A60").Value = "Programmierung"
$ws.Range("A61").Value = "- Logik"
$ws.Range("C61").Value = "2"
$ws.Range("A62").Value = "- Design "
$ws.Range("A63").Value = "Recherche"

# --- Row heights for the new header/category rows (matches the 19pt
#     rows used throughout the rest of the sheet for this font).
$ws.Rows.Item(57).RowHeight = 19
$ws.Rows.Item(58).RowHeight = 19
$ws.Rows.Item(59).RowHeight = 19
$ws.Rows.Item(60).RowHeight = 19
$ws.Rows.Item(61).RowHeight = 19
$ws.Rows.Item(62).RowHeight = 19
$ws.Rows.Item(63).RowHeight = 19

# --- Move the view/selection to the new bottom of the sheet.
$ws.Range("C64").Select()
